# GFG. Leaders in an array
# Append a new row (row 41) to the question-tracker sheet for the
# "Leaders in an array" GFG question, following the same layout used by
# the other rows (Question type in col A, Question name in col B,
# Language in col C, Date in col D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A - "GFG" tag, centered (matches the style used by other GFG rows)
$a41 = $ws.Cells.Item(41, 1)
$a41.Value = "GFG"
$a41.HorizontalAlignment = -4108

# Column B - question title, left aligned, no wrap
$b41 = $ws.Cells.Item(41, 2)
$b41.Value = "Leaders in an array"
$b41.HorizontalAlignment = -4131
$b41.WrapText = $false

# Column C - "Java" language tag, top aligned
$c41 = $ws.Cells.Item(41, 3)
$c41.Value = "Java"
$c41.VerticalAlignment = -4160

# Column D - date completed
$d41 = $ws.Cells.Item(41, 4)
$d41.NumberFormat = "d-mmm-yy"
$d41.Value = 45022

# Move / update the current selection like the saved workbook shows
$ws.Range("D42").Select()
